$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The two data rows (2 and 3) swap their Fecha (D), Volumen (M), Precio máximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values.

# Row 2 -> new values (previously held by row 3)
$ws.Range("D2").Value = 44672
$ws.Range("M2").Value = 8
$ws.Range("O2").Value = 180000
$ws.Range("P2").Value = 180000
$ws.Range("S2").Value = 180000

# Row 3 -> new values (previously held by row 2)
$ws.Range("D3").Value = 44253
$ws.Range("M3").Value = 12
$ws.Range("O3").Value = 200000
$ws.Range("P3").Value = 190000
$ws.Range("S3").Value = 190000
